$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AP1").Value = "EXTREME"
$ws.Range("AC1").Value = "HIGH CLIM"
$ws.Range("P1").Value = "REFERENCE"

$ws.Range("K18").Formula = "=SUM(K4:K9)"
$ws.Range("X18").Formula = "=SUM(X4:X9)"
$ws.Range("AK18").Formula = "=SUM(AK4:AK9)"
$ws.Range("AX18").Formula = "=SUM(AX4:AX9)"

$ws.Range("K20").Formula = "=0.1*K18*240000"
$ws.Range("X20").Formula = "=0.1*X18*240000"
$ws.Range("AK20").Formula = "=0.1*AK18*240000"
$ws.Range("AX20").Formula = "=0.1*AX18*240000"

$ws.Range("K21").Formula = "=K20*326000"
$ws.Range("L21").Formula = "=K21/100/365"
$ws.Range("X21").Formula = "=X20*326000"
$ws.Range("Y21").Formula = "=X21/100/365"
$ws.Range("AK21").Formula = "=AK20*326000"
$ws.Range("AL21").Formula = "=AK21/100/365"
$ws.Range("AX21").Formula = "=AX20*326000"
$ws.Range("AY21").Formula = "=AX21/100/365"

$ws.Range("K22").Formula = "=K21*0.005"
$ws.Range("X22").Formula = "=X21*0.005"
$ws.Range("AK22").Formula = "=AK21*0.005"
$ws.Range("AX22").Formula = "=AX21*0.005"

Write-Host "done"
